$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'30.235.85"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  +2.20%  "

$ws.Range("D3").Value = "'1.895.19"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  -0.54%  "

$ws.Range("D4").Value = "'1.001"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  -0.07%  "

$ws.Range("D5").Value = "'325.18"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +3.33%  "

$ws.Range("D6").Value = "'1.000"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -0.06%  "

$ws.Range("D7").Value = "'0.5176"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +0.38%  "

$ws.Range("D8").Value = "'0.4003"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +0.65%  "

$ws.Range("D9").Value = "'0.08420"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -0.71%  "

$ws.Range("D10").Value = "'42.64"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +0.40%  "

$ws.Range("D11").Value = "'1.114"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -0.44%  "

$ws.Range("D12").Value = "'23.36"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +12.37%  "

$ws.Range("D13").Value = "'6.439"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +2.03%  "

$ws.Range("D14").Value = "'1.892.19"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -0.32%  "

$ws.Range("D15").Value = "'7.328"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +0.02%  "

$ws.Range("E16").Value = "  -0.04%  "

$ws.Range("D17").Value = "'94.31"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +1.21%  "

$ws.Range("D18").Value = "'0.00001110"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -0.03%  "

$ws.Range("D19").Value = "'0.06644"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -1.53%  "

$ws.Range("D20").Value = "'18.22"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +1.46%  "

$ws.Range("E21").Value = "  -0.03%  "

$ws.Range("D22").Value = "'5.947"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -1.41%  "

$ws.Range("D23").Value = "'30.226.56"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +2.18%  "

$ws.Range("D24").Value = "'11.29"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +0.99%  "

$ws.Range("D25").Value = "'2.233"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +0.89%  "

$ws.Range("D26").Value = "'2.109.06"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -0.35%  "

$ws.Range("D27").Value = "'21.57"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +3.07%  "

$ws.Range("D28").Value = "'161.60"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +1.52%  "

$ws.Range("D29").Value = "'2.355"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -3.50%  "

$ws.Range("D30").Value = "'129.23"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +0.80%  "

$ws.Range("D31").Value = "'1.090"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +2.06%  "

$ws.Range("D32").Value = "'0.1054"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +0.09%  "

$ws.Range("D33").Value = "'6.087"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -1.37%  "

$ws.Range("D34").Value = "'3.748"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +2.71%  "

$ws.Range("D35").Value = "'0.02499"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +0.28%  "

$ws.Range("D36").Value = "'0.06543"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -1.03%  "

$ws.Range("D37").Value = "'5.261"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +0.42%  "

$ws.Range("D39").Value = "'1.219"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -1.18%  "

$ws.Range("D40").Value = "'11.77"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +4.07%  "

$ws.Range("B41").Value = "FraxShare"
$ws.Range("C41").Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$ws.Range("D41").Value = "'8.738"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -3.62%  "

$ws.Range("B42").Value = "TheSandbox"
$ws.Range("C42").Value = "https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand"
$ws.Range("D42").Value = "'0.6502"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -0.52%  "

$ws.Range("D43").Value = "'1.228"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -0.61%  "

$ws.Range("D44").Value = "'0.6099"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +0.11%  "

$ws.Range("D45").Value = "'13.28"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +0.86%  "

$ws.Range("E46").Value = "  +0.72%  "

$ws.Range("E47").Value = "  -0.34%  "

$ws.Range("D48").Value = "'1.236"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +0.23%  "

$ws.Range("D49").Value = "'124.60"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +0.48%  "

$ws.Range("E50").Value = "  +0.34%  "

$ws.Range("E51").Value = "  +1.21%  "
